$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($sheet, $addr, $val)
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "66.129.38"
Set-TextCell $ws "E2" "  -4.14%  "
Set-TextCell $ws "D3" "3.544.70"
Set-TextCell $ws "E3" "  -4.42%  "
Set-TextCell $ws "E4" "  +0.22%  "
Set-TextCell $ws "D5" "577.43"
Set-TextCell $ws "E5" "  -6.03%  "
Set-TextCell $ws "D6" "187.17"
Set-TextCell $ws "E6" "  -1.94%  "
Set-TextCell $ws "D7" "3.539.47"
Set-TextCell $ws "E7" "  -4.47%  "
Set-TextCell $ws "D8" "0.611"
Set-TextCell $ws "E8" "  -3.70%  "
Set-TextCell $ws "E9" "  +0.09%  "
Set-TextCell $ws "D10" "0.663"
Set-TextCell $ws "E10" "  -6.91%  "
Set-TextCell $ws "D11" "0.144"
Set-TextCell $ws "E11" "  -9.58%  "
Set-TextCell $ws "D12" "52.61"
Set-TextCell $ws "E12" "  -6.12%  "
Set-TextCell $ws "D13" "0.0000256"
Set-TextCell $ws "E13" "  -10.96%  "
Set-TextCell $ws "D14" "9.72"
Set-TextCell $ws "E14" "  -7.46%  "
Set-TextCell $ws "D15" "4.124.00"
Set-TextCell $ws "E15" "  -4.05%  "
Set-TextCell $ws "D16" "3.561.95"
Set-TextCell $ws "E16" "  -3.85%  "
Set-TextCell $ws "E17" "  -0.99%  "
Set-TextCell $ws "D18" "18.17"
Set-TextCell $ws "E18" "  -5.55%  "
Set-TextCell $ws "D19" "66.137.87"
Set-TextCell $ws "E19" "  -3.81%  "
Set-TextCell $ws "D20" "12.07"
Set-TextCell $ws "E20" "  -6.27%  "
Set-TextCell $ws "E21" "  -7.19%  "
Set-TextCell $ws "D22" "390.83"
Set-TextCell $ws "E22" "  -4.60%  "
Set-TextCell $ws "D23" "4.28"
Set-TextCell $ws "E23" "  -6.84%  "
Set-TextCell $ws "D24" "85.16"
Set-TextCell $ws "E24" "  -4.21%  "
Set-TextCell $ws "D25" "10.98"
Set-TextCell $ws "E25" "  +1.03%  "
Set-TextCell $ws "E26" "  -4.98%  "
Set-TextCell $ws "D27" "12.31"
Set-TextCell $ws "E27" "  -3.29%  "
Set-TextCell $ws "E28" "  +0.03%  "
Set-TextCell $ws "D29" "3.50"
Set-TextCell $ws "E29" "  -6.43%  "
Set-TextCell $ws "D30" "8.83"
Set-TextCell $ws "E30" "  -8.16%  "
Set-TextCell $ws "D31" "30.79"
Set-TextCell $ws "E31" "  -6.58%  "
Set-TextCell $ws "D32" "7.07"
Set-TextCell $ws "E32" "  -2.03%  "
Set-TextCell $ws "D33" "627.67"
Set-TextCell $ws "E33" "  +0.67%  "
Set-TextCell $ws "D34" "12.08"
Set-TextCell $ws "E34" "  -4.03%  "
Set-TextCell $ws "D35" "63.35"
Set-TextCell $ws "E35" "  -3.47%  "
Set-TextCell $ws "D36" "0.112"
Set-TextCell $ws "E36" "  -7.33%  "
Set-TextCell $ws "D37" "41.06"
Set-TextCell $ws "E37" "  -7.37%  "
Set-TextCell $ws "E38" "  +0.15%  "
Set-TextCell $ws "D39" "0.393"
Set-TextCell $ws "E39" "  -4.16%  "
Set-TextCell $ws "D40" "0.0₃0758"
Set-TextCell $ws "E40" "  -6.09%  "
Set-TextCell $ws "B41" "FirstDigitalUSD"
Set-TextCell $ws "C41" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D41" "1.00"
Set-TextCell $ws "E41" "  -0.05%  "
Set-TextCell $ws "B42" "Kaspa"
Set-TextCell $ws "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws "D42" "0.130"
Set-TextCell $ws "E42" "  -6.79%  "
Set-TextCell $ws "D43" "2.970.62"
Set-TextCell $ws "E43" "  +3.82%  "
Set-TextCell $ws "D44" "2.78"
Set-TextCell $ws "E44" "  -7.64%  "
Set-TextCell $ws "D45" "2.47"
Set-TextCell $ws "E45" "  -4.77%  "
Set-TextCell $ws "D46" "0.0404"
Set-TextCell $ws "E46" "  -8.31%  "
Set-TextCell $ws "D47" "0.129"
Set-TextCell $ws "E47" "  -7.13%  "
Set-TextCell $ws "D48" "3.06"
Set-TextCell $ws "E48" "  -2.02%  "
Set-TextCell $ws "D49" "137.81"
Set-TextCell $ws "E49" "  -2.56%  "
Set-TextCell $ws "D50" "8.38"
Set-TextCell $ws "E50" "  -7.43%  "
Set-TextCell $ws "B51" "WEMIXToken"
Set-TextCell $ws "C51" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D51" "2.49"
Set-TextCell $ws "E51" "  -8.55%  "
